$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 3 for the columns that differ
# (D: Fecha, M: Volumen, O: Precio máximo, P: Precio promedio ponderado, S: Precio $/Kg)

$ws.Range("D2").Value = 44672
$ws.Range("M2").Value = 8
$ws.Range("O2").Value = 180000
$ws.Range("P2").Value = 180000
$ws.Range("S2").Value = 180000

$ws.Range("D3").Value = 44993
$ws.Range("M3").Value = 14
$ws.Range("O3").Value = 200000
$ws.Range("P3").Value = 190000
$ws.Range("S3").Value = 190000
